# wizard_of_oz_experiment_data_solution.xlsx
#
# "fixed problem with database (error was in config excel)"
#
# The "systemId" column in the solutionTickets table was driven by a
# calculated formula (IFERROR(<prev>+1,0)) which kept breaking the
# downstream database import. The fix freezes the column to its current
# static values (0..6) and fills in the missing
# "ticketDescriptionHighlighting" value for the newest ticket row (row 8,
# column G) with the same empty-highlight default ("[]") used by the
# other rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solutionTickets")

# --- Fill in the missing ticketDescriptionHighlighting value for the
#     last ticket row, matching the default used elsewhere in the column.
$ws.Range("G8").Value = "[]"

# --- Freeze the "systemId" helper column: replace the auto-incrementing
#     formulas in A2:A8 with their already-computed static values so the
#     sheet no longer depends on the flaky formula chain.
$col = $ws.Range("A2:A8")
$col.Value2 = $col.Value2

# --- Reflect where the user ended up after making the edit.
$ws.Range("G8").Select() | Out-Null
